$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume/change (E) columns
# Values that look like plain numbers are apostrophe-prefixed so Excel
# keeps storing them as text (matching the original inline-string cells)
# instead of silently converting them to numeric values.

$ws.Range("D2").Value = '26.800.25'
$ws.Range("E2").Value = '  +1.04%  '
$ws.Range("D3").Value = '1.649.36'
$ws.Range("E3").Value = '  +1.33%  '
$ws.Range("E4").Value = '  +0.49%  '
$ws.Range("D5").Value = '''216.64'
$ws.Range("E5").Value = '  +1.67%  '
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("E7").Value = '  +0.71%  '
$ws.Range("E9").Value = '  +0.29%  '
$ws.Range("D10").Value = '''19.14'
$ws.Range("E10").Value = '  +1.50%  '
$ws.Range("D12").Value = '1.879.44'
$ws.Range("D13").Value = '1.650.94'
$ws.Range("E13").Value = '  +1.60%  '
$ws.Range("D14").Value = '''4.19'
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("E15").Value = '  +1.72%  '
$ws.Range("D16").Value = '''65.34'
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("D17").Value = '26.809.60'
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("D19").Value = '''218.22'
$ws.Range("E19").Value = '  +1.72%  '
$ws.Range("E20").Value = '  +0.62%  '
$ws.Range("D21").Value = '''4.36'
$ws.Range("E21").Value = '  +1.62%  '
$ws.Range("D22").Value = '''2.45'
$ws.Range("E22").Value = '  +18.25%  '
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("E24").Value = '  +1.91%  '
$ws.Range("D25").Value = '''146.66'
$ws.Range("E25").Value = '  -1.00%  '
$ws.Range("E26").Value = '  +0.77%  '
$ws.Range("E27").Value = '  -0.40%  '
$ws.Range("E28").Value = '  +3.83%  '
$ws.Range("D29").Value = '''15.74'
$ws.Range("E29").Value = '  +1.29%  '
$ws.Range("D30").Value = '''0.0519'
$ws.Range("E30").Value = '  +1.65%  '
$ws.Range("E31").Value = '  +1.71%  '
$ws.Range("D32").Value = '''3.35'
$ws.Range("E32").Value = '  +0.22%  '
$ws.Range("E33").Value = '  +1.23%  '
$ws.Range("D34").Value = '1.283.44'
$ws.Range("E34").Value = '  +3.45%  '
$ws.Range("E35").Value = '  +2.99%  '
$ws.Range("E36").Value = '  +3.17%  '
$ws.Range("E37").Value = '  +2.33%  '
$ws.Range("E38").Value = '  +5.80%  '
$ws.Range("D39").Value = '''0.828'
$ws.Range("E39").Value = '  +4.30%  '
$ws.Range("E40").Value = '  +0.85%  '
$ws.Range("D41").Value = '''0.814'
$ws.Range("E41").Value = '  +1.98%  '
$ws.Range("D42").Value = '''2.25'
$ws.Range("E42").Value = '  -1.32%  '
$ws.Range("D43").Value = '''5.47'
$ws.Range("E43").Value = '  +2.50%  '
$ws.Range("D44").Value = '1.789.63'
$ws.Range("E44").Value = '  +1.44%  '
$ws.Range("D45").Value = '''92.06'
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("D46").Value = '''59.61'
$ws.Range("E46").Value = '  +8.51%  '
$ws.Range("E47").Value = '  +1.38%  '
$ws.Range("D48").Value = '''0.0515'
$ws.Range("E48").Value = '  +1.24%  '
$ws.Range("E49").Value = '  +3.45%  '
$ws.Range("D50").Value = '''0.0970'
$ws.Range("E50").Value = '  +1.64%  '
$ws.Range("D51").Value = '''0.408'
$ws.Range("E51").Value = '  +0.43%  '
